$wb = $excel.ActiveWorkbook

# --- Sheet ALC (index 1) ---
$ws = $wb.Worksheets.Item(1)
# Row 62
$ws.Range("H62").Value = 4833
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 5291.25
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 5291.25
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -6539.25

# Row 65
$ws.Range("H65").Value = 4833
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 5291.25
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 26456.25
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -32696.25

# Row 97
$ws.Range("H97").Value = 1072.3
$ws.Range("J97").Value = 1072.3
$ws.Range("L97").Value = 3216.9
$ws.Range("N97").Value = -4208.9

# Row 112
$ws.Range("H112").Value = 4546642
$ws.Range("I112").Value = 2798
$ws.Range("J112").Value = 5001026.5
$ws.Range("K112").Value = 8394
$ws.Range("L112").Value = 15003079.5
$ws.Range("M112").Value = -7286
$ws.Range("N112").Value = -15005295.5

# Row 129
$ws.Range("H129").Value = 5103670.5
$ws.Range("J129").Value = 1798
$ws.Range("L129").Value = 5394
$ws.Range("N129").Value = -15394

# Row 132
$ws.Range("H132").Value = 3637734
$ws.Range("I132").Value = 4082534.8
$ws.Range("J132").Value = 5194.3335
$ws.Range("K132").Value = 12247604.4
$ws.Range("L132").Value = 15583.0005
$ws.Range("M132").Value = -12245074.4
$ws.Range("N132").Value = -20643.0005

# Row 138
$ws.Range("H138").Value = 3483.976
$ws.Range("I138").Value = 2063.1
$ws.Range("J138").Value = 7036.1665
$ws.Range("K138").Value = 6189.299999999999
$ws.Range("L138").Value = 21108.4995
$ws.Range("M138").Value = -1049.299999999999
$ws.Range("N138").Value = -31388.4995

# --- Sheet ARM (index 2) ---
$ws = $wb.Worksheets.Item(2)
# Row 23
$ws.Range("H23").Value = 26111.777
$ws.Range("I23").Value = 44003
$ws.Range("J23").Value = 21000
$ws.Range("K23").Value = 44003
$ws.Range("L23").Value = 21000
$ws.Range("M23").Value = -43744
$ws.Range("N23").Value = -21518

# Row 58
$ws.Range("H58").Value = 30000
$ws.Range("J58").Value = 30000
$ws.Range("L58").Value = 30000
$ws.Range("N58").Value = -30860

# Row 74
$ws.Range("H74").Value = 790.3125
$ws.Range("I74").Value = 699.62964
$ws.Range("J74").Value = 1280
$ws.Range("K74").Value = 699.62964
$ws.Range("L74").Value = 1280
$ws.Range("M74").Value = 174.37036
$ws.Range("N74").Value = -3028

# Row 77
$ws.Range("H77").Value = 790.3125
$ws.Range("I77").Value = 699.62964
$ws.Range("J77").Value = 1280
$ws.Range("K77").Value = 3498.1482
$ws.Range("L77").Value = 6400
$ws.Range("M77").Value = 869.8517999999999
$ws.Range("N77").Value = -15136

# Row 122
$ws.Range("H122").Value = 4173.5454
$ws.Range("I122").Value = 1982
$ws.Range("J122").Value = 5999.8335
$ws.Range("K122").Value = 5946
$ws.Range("L122").Value = 17999.5005
$ws.Range("M122").Value = -3496
$ws.Range("N122").Value = -22899.5005

# --- Sheet CRP (index 4) ---
$ws = $wb.Worksheets.Item(4)
# Row 99
$ws.Range("H99").Value = 4666.6665
$ws.Range("J99").Value = 4666.6665
$ws.Range("L99").Value = 4666.6665
$ws.Range("N99").Value = -7662.6665

# Row 105
$ws.Range("H105").Value = 1928.5714
$ws.Range("I105").Value = 1915.75
$ws.Range("J105").Value = 1945.6666
$ws.Range("K105").Value = 1915.75
$ws.Range("L105").Value = 1945.6666
$ws.Range("M105").Value = -168.75
$ws.Range("N105").Value = -5439.6666

# Row 126
$ws.Range("H126").Value = 4666.6665
$ws.Range("J126").Value = 4666.6665
$ws.Range("L126").Value = 13999.9995
$ws.Range("N126").Value = -18939.9995

# Row 132
$ws.Range("H132").Value = 4464.5293
$ws.Range("I132").Value = 3111
$ws.Range("K132").Value = 9333
$ws.Range("M132").Value = -6803

# --- Sheet CUL (index 5) ---
$ws = $wb.Worksheets.Item(5)
# Row 68
$ws.Range("H68").Value = 2152.0178
$ws.Range("I68").Value = 705.25
$ws.Range("J68").Value = 2955.7778
$ws.Range("K68").Value = 2115.75
$ws.Range("L68").Value = 8867.3334
$ws.Range("M68").Value = -1304.75
$ws.Range("N68").Value = -10489.3334

# Row 71
$ws.Range("H71").Value = 2152.0178
$ws.Range("I71").Value = 705.25
$ws.Range("J71").Value = 2955.7778
$ws.Range("K71").Value = 6347.25
$ws.Range("L71").Value = 26602.0002
$ws.Range("M71").Value = -2291.25
$ws.Range("N71").Value = -34714.00019999999

# Row 131
$ws.Range("H131").Value = 2195.5264
$ws.Range("J131").Value = 1286.25
$ws.Range("L131").Value = 3858.75
$ws.Range("N131").Value = -13938.75

# Row 132
$ws.Range("H132").Value = 4672
$ws.Range("I132").Value = 1901.3334
$ws.Range("J132").Value = 6750
$ws.Range("K132").Value = 17112.0006
$ws.Range("L132").Value = 60750
$ws.Range("M132").Value = -14582.0006
$ws.Range("N132").Value = -65810

# --- Sheet GSM (index 6) ---
$ws = $wb.Worksheets.Item(6)
# Row 52
$ws.Range("H52").Value = 44000
$ws.Range("J52").Value = 44000
$ws.Range("L52").Value = 44000
$ws.Range("N52").Value = -44518

# Row 107
$ws.Range("H107").Value = 656.3913
$ws.Range("I107").Value = 247
$ws.Range("J107").Value = 919.5714
$ws.Range("K107").Value = 247
$ws.Range("L107").Value = 919.5714
$ws.Range("M107").Value = 1673
$ws.Range("N107").Value = -4759.5714

# --- Sheet LTW (index 7) ---
$ws = $wb.Worksheets.Item(7)
# Row 7
$ws.Range("H7").Value = 2516.7856
$ws.Range("I7").Value = 1797.5
$ws.Range("K7").Value = 1797.5
$ws.Range("M7").Value = -1685.5

# Row 22
$ws.Range("H22").Value = 125001810
$ws.Range("I22").Value = 142857780
$ws.Range("J22").Value = 10000
$ws.Range("K22").Value = 142857780
$ws.Range("L22").Value = 10000
$ws.Range("M22").Value = -142857485
$ws.Range("N22").Value = -10590

# Row 27
$ws.Range("H27").Value = 125001810
$ws.Range("I27").Value = 142857780
$ws.Range("J27").Value = 10000
$ws.Range("K27").Value = 142857780
$ws.Range("L27").Value = 10000
$ws.Range("M27").Value = -142857673
$ws.Range("N27").Value = -10214

# Row 93
$ws.Range("H93").Value = 5244.5713
$ws.Range("I93").Value = 3942.4
$ws.Range("K93").Value = 3942.4
$ws.Range("M93").Value = -2694.4

# Row 126
$ws.Range("H126").Value = 2516.7856
$ws.Range("I126").Value = 1797.5
$ws.Range("K126").Value = 5392.5
$ws.Range("M126").Value = -2922.5

# Row 132
$ws.Range("H132").Value = 3163.0386
$ws.Range("I132").Value = 2219.3845
$ws.Range("J132").Value = 4106.6924
$ws.Range("K132").Value = 6658.1535
$ws.Range("L132").Value = 12320.0772
$ws.Range("M132").Value = -4128.1535
$ws.Range("N132").Value = -17380.0772

# --- Sheet WVR (index 8) ---
$ws = $wb.Worksheets.Item(8)
# Row 126
$ws.Range("H126").Value = 3032325.5
$ws.Range("I126").Value = 1664.8422
$ws.Range("J126").Value = 7145365
$ws.Range("K126").Value = 4994.5266
$ws.Range("L126").Value = 21436095
$ws.Range("M126").Value = -2524.5266
$ws.Range("N126").Value = -21441035

Write-Output "All 29 rows updated across 7 sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR)."
